$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last currently-used data row (before the edit) is row 21 (header is row 1).
$lastOldRow = 21
$shift = 2

# Shift the existing data rows (2..21) down by two rows (to 4..23) by
# copying values directly, starting from the bottom so we don't clobber
# rows before they are read. This avoids Excel's row-insert behavior of
# copying cell formatting from adjoining rows.
for ($r = $lastOldRow; $r -ge 2; $r--) {
    $destRow = $r + $shift
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value()
}

# Fill in the two newly opened rows (2 and 3) with their new values.
$newTopRows = @(
    @(0.04686117172241161, -0.03032520040869657, 0.02842492796480637),
    @(0.01706874370574948, 0.07629761099815366, -0.1732735317200423)
)

for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTopRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTopRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTopRows[$i][2]
}

# Append eight new rows of data at the bottom (new rows 24-31).
$newBottomRows = @(
    @(-0.4282075166702048, -0.7951091900468019, -1.767483308911331),
    @(-2.543609619140639, 0.4462372660637008, -1.760738492012012),
    @(0.2467263936996389, -0.1943315342068692, -2.037849001586441),
    @(0.9171624183654843, 0.3103487230837345, 0.1225722581148094),
    @(-0.926007807254792, 0.2953229788690807, 2.136403992772098),
    @(-0.3645055294036889, 0.3592699170112605, 1.912454850971703),
    @(0.7052955031394972, 0.3671576231718077, -0.5781752001494165),
    @(-0.5357744693756123, -0.6449819654226362, -0.8872665241360714)
)

$startRow = 24
for ($i = 0; $i -lt $newBottomRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newBottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newBottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newBottomRows[$i][2]
}
